# Fix the tutorial link in Lab1Instructions.
#
# 1. Remove the colon after "Do the first five parts of" (the trailing
#    space is kept, but split into its own run further below so the
#    structure matches a real Word hyperlink-edit session).
# 2. Point the hyperlink at the new ASP.NET MVC tutorial URL and update
#    its visible text to match.
# 3. The editing session leaves the "_GoBack" bookmark positioned right
#    after the hyperlink (where the cursor was after editing it) instead
#    of at the end of the paragraph.

$d = $word.ActiveDocument

# --- Step 1: drop the colon ("Do the first five parts of: " -> "Do the first five parts of ") ---
$scope = $d.Range(0, $d.Content.End)
$scope.Find.Execute("Do the first five parts of: ", $true, $false, $false, $false, $false, $true, 1, $false, "Do the first five parts of ", 2)

# --- Step 2: split the trailing space into its own run (mirrors the run
#     boundary Word leaves behind after editing the hyperlink text) ---
$scope2 = $d.Range(0, $d.Content.End)
$scope2.Find.Execute("Do the first five parts of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spaceRng = $d.Range($scope2.End - 1, $scope2.End)
$spaceRng.Bold = 1
$spaceRng.Bold = 0

# --- Step 3: repoint the hyperlink at the new tutorial page ---
$h = $d.Hyperlinks(1)
$h.TextToDisplay = "http://www.asp.net/mvc/overview/getting-started/introduction/getting-started"
$h.Address = "http://www.asp.net/mvc/overview/getting-started/introduction/getting-started"

# --- Step 4: move the "_GoBack" bookmark to right after the hyperlink ---
$h2 = $d.Hyperlinks(1)
$afterLink = $h2.Range.End
$d.Bookmarks.Add("_GoBack", $d.Range($afterLink, $afterLink))
